# Updates cryptos list price (D) and 1h volume change (E) columns
# to match the latest scrape, per the Jan 8 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.004.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "'2.267.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  -0.65%  "
$ws.Range("D5").Value = "'302.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("D6").Value = "'95.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("E7").Value = "  -1.63%  "
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("D9").Value = "'0.509"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.37%  "
$ws.Range("D10").Value = "'34.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.15%  "
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("E12").Value = "  -0.61%  "
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").Value = "'2.615.62"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("D15").Value = "'2.271.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("D16").Value = "'13.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "'0.799"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.23%  "
$ws.Range("D18").Value = "'44.912.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("D19").Value = "'12.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.85%  "
$ws.Range("D20").Value = "'0.0₃0923"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.03%  "
$ws.Range("D21").Value = "'6.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.34%  "
$ws.Range("D22").Value = "'65.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'238.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("E24").Value = "  -1.98%  "
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("E26").Value = "  -4.79%  "
$ws.Range("D27").Value = "'41.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.73%  "
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("E29").Value = "  -2.99%  "
$ws.Range("D30").Value = "'19.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.08%  "
$ws.Range("D31").Value = "'153.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("D32").Value = "'5.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.27%  "
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("E34").Value = "  -2.45%  "
$ws.Range("D35").Value = "'2.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.39%  "
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("D37").Value = "'0.104"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.09%  "
$ws.Range("D38").Value = "'1.77"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.02%  "
$ws.Range("D39").Value = "'4.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.92%  "
$ws.Range("D40").Value = "'0.0313"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.33%  "
$ws.Range("D41").Value = "'3.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.15%  "
$ws.Range("D42").Value = "'13.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.64%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("E44").Value = "  +12.19%  "
$ws.Range("D45").Value = "'1.751.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.49%  "
$ws.Range("D46").Value = "'0.198"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.61%  "
$ws.Range("D47").Value = "'70.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.89%  "
$ws.Range("D48").Value = "'75.95"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.34%  "
$ws.Range("D49").Value = "'96.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.09%  "
$ws.Range("D50").Value = "'53.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("D51").Value = "'4.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.81%  "
